$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows (56-68) to append under the existing results table.
# Column C: either a numeric id or a text label (shared-string-backed).
# Column F: numeric measurement.
# Column G: numeric count.
$rows = @(
    @{ Row = 56; C = 45096;    F = 53.078;    G = 2 },
    @{ Row = 57; C = 188091;   F = 24.0496;   G = 3 },
    @{ Row = 58; C = "test";   F = 38.9665;   G = 2 },
    @{ Row = 59; C = 253036;   F = 28.33;     G = 2 },
    @{ Row = 60; C = 42049;    F = 24.6923;   G = 3 },
    @{ Row = 61; C = 35070;    F = 17.7065;   G = 5 },
    @{ Row = 62; C = 163014;   F = 18.2836;   G = 4 },
    @{ Row = 63; C = 124084;   F = 51.0747;   G = 2 },
    @{ Row = 64; C = 296059;   F = 19.5217;   G = 4 },
    @{ Row = 65; C = 176035;   F = 39.8875;   G = 2 },
    @{ Row = 66; C = 295087;   F = 20.5204;   G = 5 },
    @{ Row = 67; C = 216066;   F = 16.6339;   G = 6 },
    @{ Row = 68; C = "41004-2"; F = 36.7055;  G = 2 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 6).Value = $r.F
    $ws.Cells.Item($r.Row, 7).Value = $r.G
}

# Match the saved view state from the diff (scrolled so row 58 / col B is the
# top-left visible cell, with F69 as the active selection).
$excel.ActiveWindow.ScrollRow = 58
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("F69").Select()
